$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the "in_service" boolean column (E) to TRUE for rows 10 through 15
$ws.Range("E10:E15").Value = $true
